{"js": "const body = context.document.body;\nconst replacements = [\n  [\"176\u00d72=\", \"813\u00d77=\"],\n  [\"135\u00d73=\", \"161\u00d77=\"],\n  [\"954\u00d72=\", \"558\u00d79=\"],\n  [\"713\u00d75=\", \"595\u00d78=\"],\n  [\"409\u00d79=\", \"283\u00d73=\"],\n  [\"132\u00d79=\", \"102\u00d75=\"],\n  [\"643\u00d73=\", \"579\u00d79=\"],\n  [\"516\u00d74=\", \"455\u00d75=\"],\n  [\"166\u00d75=\", \"819\u00d73=\"],\n  [\"828\u00d74=\", \"784\u00d79=\"],\n  [\"254\u00d76=\", \"676\u00d75=\"],\n  [\"367\u00d73=\", \"480\u00d79=\"],\n  [\"254\u00d77=\", \"206\u00d72=\"],\n  [\"373\u00d79=\", \"405\u00d75=\"],\n  [\"442\u00d72=\", \"549\u00d73=\"],\n  [\"221\u00d73=\", \"426\u00d75=\"],\n  [\"903\u00d76=\", \"805\u00d74=\"],\n  [\"326\u00d78=\", \"568\u00d73=\"],\n  [\"976\u00d73=\", \"536\u00d72=\"],\n  [\"385\u00d73=\", \"874\u00d74=\"],\n  [\"453\u00d79=\", \"554\u00d72=\"],\n  [\"231\u00d77=\", \"344\u00d73=\"],\n  [\"937\u00d79=\", \"405\u00d73=\"],\n  [\"602\u00d76=\", \"976\u00d77=\"],\n  [\"624\u00d77=\", \"112\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"176\u00d72=\"; New = \"813\u00d77=\" },\n    @{ Old = \"135\u00d73=\"; New = \"161\u00d77=\" },\n    @{ Old = \"954\u00d72=\"; New = \"558\u00d79=\" },\n    @{ Old = \"713\u00d75=\"; New = \"595\u00d78=\" },\n    @{ Old = \"409\u00d79=\"; New = \"283\u00d73=\" },\n    @{ Old = \"132\u00d79=\"; New = \"102\u00d75=\" },\n    @{ Old = \"643\u00d73=\"; New = \"579\u00d79=\" },\n    @{ Old = \"516\u00d74=\"; New = \"455\u00d75=\" },\n    @{ Old = \"166\u00d75=\"; New = \"819\u00d73=\" },\n    @{ Old = \"828\u00d74=\"; New = \"784\u00d79=\" },\n    @{ Old = \"254\u00d76=\"; New = \"676\u00d75=\" },\n    @{ Old = \"367\u00d73=\"; New = \"480\u00d79=\" },\n    @{ Old = \"254\u00d77=\"; New = \"206\u00d72=\" },\n    @{ Old = \"373\u00d79=\"; New = \"405\u00d75=\" },\n    @{ Old = \"442\u00d72=\"; New = \"549\u00d73=\" },\n    @{ Old = \"221\u00d73=\"; New = \"426\u00d75=\" },\n    @{ Old = \"903\u00d76=\"; New = \"805\u00d74=\" },\n    @{ Old = \"326\u00d78=\"; New = \"568\u00d73=\" },\n    @{ Old = \"976\u00d73=\"; New = \"536\u00d72=\" },\n    @{ Old = \"385\u00d73=\"; New = \"874\u00d74=\" },\n    @{ Old = \"453\u00d79=\"; New = \"554\u00d72=\" },\n    @{ Old = \"231\u00d77=\"; New = \"344\u00d73=\" },\n    @{ Old = \"937\u00d79=\"; New = \"405\u00d73=\" },\n    @{ Old = \"602\u00d76=\"; New = \"976\u00d77=\" },\n    @{ Old = \"624\u00d77=\"; New = \"112\u00d74=\" },\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute([ref]$find.Text, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]2) | Out-Null\n}\n"}
